# Update United States 2010 (column G) and a couple of United States 2020
# (column H) figures in the demographics table on Sheet1.
#
# All of the values in this table are stored as text (not numbers), so a
# leading apostrophe is used to force each assignment to remain text
# instead of being auto-converted to a numeric value by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = "'0.48"
$ws.Range("G6").Value  = "'0.04"
$ws.Range("G10").Value = "'0.6"
$ws.Range("H10").Value = "'0.59"
$ws.Range("G11").Value = "'0.19"
$ws.Range("G13").Value = "'1.99"
$ws.Range("G15").Value = "'0.19"
$ws.Range("H15").Value = "'0.21"
$ws.Range("G17").Value = "'1111430"
$ws.Range("G19").Value = "'0.54"
$ws.Range("G20").Value = "'0.3"
$ws.Range("G25").Value = "'0.13"
$ws.Range("G26").Value = "'0.5"
$ws.Range("H26").Value = "'0.53"
$ws.Range("G29").Value = "'2.15"
$ws.Range("G30").Value = "'0.21"
$ws.Range("G31").Value = "'0.17"
$ws.Range("H31").Value = "'0.18"
$ws.Range("G32").Value = "'0.72"
$ws.Range("G33").Value = "'889000"
